$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45951, 7441.5535039366, 7105.58627936105, 9404, 3261.88802402101, 40.1447626409193),
    @(45952, 7237.61420618573, 6814.09503889483, 3620, 3325.93590967206, 271.667956190287),
    @(45953, 6964.57551660774, 6562.8032075483, 3620, 3136.5748830955, 253.307420443492),
    @(45954, 7010.33215944304, 6188.35804127123, 3620, 3186.50811711493, 239.786089932757),
    @(45955, 3094.05598362749, 4882.69502771512, 3620, 3630.90364287561, 203.899944607947),
    @(45956, 2228.83412809025, 4391.37737199534, 3620, 2820.83705457257, 149.675601106996),
    @(45957, 7224.82758554897, 7519.12319514416, 3620, 3510.4863574464, 308.73373135794),
    @(45958, 7224.82758554897, 7893.54552457967, 3620, 3510.4863574464, 324.334661751087),
    @(45959, 7224.82758554897, 7534.68716881837, 3620, 3510.4863574464, 309.382230261032),
    @(45960, 7224.82758554897, 7348.99032276079, 3620, 3510.4863574464, 301.6448616753),
    @(45961, 6107.35191671871, 5662.73405803112, 3620, 3445.86867539698, 228.691780559504),
    @(45962, 2557.16655997904, 3750.3935641774, 6652, 3795.7240222106, 37.2548994328335),
    @(45963, 2557.16655997904, 3802.2009669056, 6652, 3795.7240222106, 39.413541213175),
    @(45964, 8571.71877057272, 6885.96182504027, 6652, 4684.40102234949, 204.931785307907)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $ws.Cells.Item($startRow + $i, $j + 1).Value = $rowVals[$j]
    }
}
